$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2777.2222
$ws.Range("J19").Value = 2749
$ws.Range("L19").Value = 2749
$ws.Range("N19").Value = -3099
$ws.Range("H41").Value = 19232498
$ws.Range("I41").Value = 664.3570999999999
$ws.Range("J41").Value = 41669640
$ws.Range("K41").Value = 664.3570999999999
$ws.Range("L41").Value = 41669640
$ws.Range("M41").Value = -224.3570999999999
$ws.Range("N41").Value = -41670520
$ws.Range("H62").Value = 62501250
$ws.Range("I62").Value = 125000000
$ws.Range("K62").Value = 125000000
$ws.Range("M62").Value = -124999376
$ws.Range("H64").Value = 5915.25
$ws.Range("I64").Value = 4997.385
$ws.Range("J64").Value = 7000
$ws.Range("K64").Value = 4997.385
$ws.Range("L64").Value = 7000
$ws.Range("M64").Value = -4749.385
$ws.Range("N64").Value = -7496
$ws.Range("H65").Value = 62501250
$ws.Range("I65").Value = 125000000
$ws.Range("K65").Value = 625000000
$ws.Range("M65").Value = -624996880
$ws.Range("H67").Value = 5915.25
$ws.Range("I67").Value = 4997.385
$ws.Range("J67").Value = 7000
$ws.Range("K67").Value = 4997.385
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = -4139.385
$ws.Range("N67").Value = -8716
$ws.Range("H80").Value = 1141.0869
$ws.Range("I80").Value = 780.26666
$ws.Range("J80").Value = 1817.625
$ws.Range("K80").Value = 2340.79998
$ws.Range("L80").Value = 5452.875
$ws.Range("M80").Value = -1342.79998
$ws.Range("N80").Value = -7448.875
$ws.Range("H83").Value = 1141.0869
$ws.Range("I83").Value = 780.26666
$ws.Range("J83").Value = 1817.625
$ws.Range("K83").Value = 7022.39994
$ws.Range("L83").Value = 16358.625
$ws.Range("M83").Value = -2030.39994
$ws.Range("N83").Value = -26342.625
$ws.Range("H103").Value = 334.875
$ws.Range("J103").Value = 334.875
$ws.Range("L103").Value = 1004.625
$ws.Range("N103").Value = -2176.625
$ws.Range("H107").Value = 60047.94
$ws.Range("I107").Value = 101457.3
$ws.Range("J107").Value = 891.7143
$ws.Range("K107").Value = 101457.3
$ws.Range("L107").Value = 891.7143
$ws.Range("M107").Value = -99537.3
$ws.Range("N107").Value = -4731.7143
$ws.Range("H132").Value = 1957.5128
$ws.Range("I132").Value = 1363.8438
$ws.Range("J132").Value = 4671.4287
$ws.Range("K132").Value = 4091.5314
$ws.Range("L132").Value = 14014.2861
$ws.Range("M132").Value = -1561.5314
$ws.Range("N132").Value = -19074.2861
$ws.Range("H137").Value = 2871.186
$ws.Range("I137").Value = 2158.3845
$ws.Range("K137").Value = 6475.1535
$ws.Range("M137").Value = -3925.1535

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4690.1
$ws.Range("I2").Value = 4287.6875
$ws.Range("J2").Value = 6299.75
$ws.Range("K2").Value = 4287.6875
$ws.Range("L2").Value = 6299.75
$ws.Range("M2").Value = -4174.6875
$ws.Range("N2").Value = -6525.75
$ws.Range("H45").Value = 9900
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 9900
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9900
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -10654
$ws.Range("H74").Value = 12118.962
$ws.Range("I74").Value = 14304.85
$ws.Range("K74").Value = 14304.85
$ws.Range("M74").Value = -13430.85
$ws.Range("H77").Value = 12118.962
$ws.Range("I77").Value = 14304.85
$ws.Range("K77").Value = 71524.25
$ws.Range("M77").Value = -67156.25
$ws.Range("H97").Value = 1713.7
$ws.Range("I97").Value = 1713.7
$ws.Range("K97").Value = 1713.7
$ws.Range("M97").Value = -1217.7
$ws.Range("H110").Value = 834051.7
$ws.Range("I110").Value = 834051.7
$ws.Range("K110").Value = 834051.7
$ws.Range("M110").Value = -832006.7
$ws.Range("H116").Value = 4690.1
$ws.Range("I116").Value = 4287.6875
$ws.Range("J116").Value = 6299.75
$ws.Range("K116").Value = 4287.6875
$ws.Range("L116").Value = 6299.75
$ws.Range("M116").Value = -1993.6875
$ws.Range("N116").Value = -10887.75
$ws.Range("H122").Value = 5028.5713
$ws.Range("I122").Value = 4110.75
$ws.Range("K122").Value = 12332.25
$ws.Range("M122").Value = -9882.25
$ws.Range("H132").Value = 8532.879999999999
$ws.Range("I132").Value = 4253.5454
$ws.Range("K132").Value = 12760.6362
$ws.Range("M132").Value = -10230.6362

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4690.1
$ws.Range("I3").Value = 4287.6875
$ws.Range("J3").Value = 6299.75
$ws.Range("K3").Value = 4287.6875
$ws.Range("L3").Value = 6299.75
$ws.Range("M3").Value = -4173.6875
$ws.Range("N3").Value = -6527.75
$ws.Range("H35").Value = 119900
$ws.Range("J35").Value = 119900
$ws.Range("L35").Value = 119900
$ws.Range("N35").Value = -120520
$ws.Range("H94").Value = 1085.3125
$ws.Range("I94").Value = 847.5
$ws.Range("K94").Value = 847.5
$ws.Range("M94").Value = -396.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3450.6875
$ws.Range("I16").Value = 2136.5454
$ws.Range("J16").Value = 6341.8
$ws.Range("K16").Value = 2136.5454
$ws.Range("L16").Value = 6341.8
$ws.Range("M16").Value = -1849.5454
$ws.Range("N16").Value = -6915.8
$ws.Range("H31").Value = 3212.6956
$ws.Range("I31").Value = 1774.5
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 1774.5
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -1479.5
$ws.Range("N31").Value = -7090
$ws.Range("H34").Value = 3212.6956
$ws.Range("I34").Value = 1774.5
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 1774.5
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -1572.5
$ws.Range("N34").Value = -6904
$ws.Range("H58").Value = 458231.38
$ws.Range("I58").Value = 834739.25
$ws.Range("K58").Value = 834739.25
$ws.Range("M58").Value = -834536.25
$ws.Range("H99").Value = 6583
$ws.Range("I99").Value = 5998.5
$ws.Range("J99").Value = 6699.9
$ws.Range("K99").Value = 5998.5
$ws.Range("L99").Value = 6699.9
$ws.Range("M99").Value = -4500.5
$ws.Range("N99").Value = -9695.9
$ws.Range("H105").Value = 1220.6666
$ws.Range("I105").Value = 1189.5454
$ws.Range("J105").Value = 1357.6
$ws.Range("K105").Value = 1189.5454
$ws.Range("L105").Value = 1357.6
$ws.Range("M105").Value = 557.4546
$ws.Range("N105").Value = -4851.6
$ws.Range("H107").Value = 2156.9614
$ws.Range("I107").Value = 1815.125
$ws.Range("J107").Value = 2703.9
$ws.Range("K107").Value = 1815.125
$ws.Range("L107").Value = 2703.9
$ws.Range("M107").Value = 104.875
$ws.Range("N107").Value = -6543.9
$ws.Range("H113").Value = 3450.6875
$ws.Range("I113").Value = 2136.5454
$ws.Range("J113").Value = 6341.8
$ws.Range("K113").Value = 2136.5454
$ws.Range("L113").Value = 6341.8
$ws.Range("M113").Value = 33.45460000000003
$ws.Range("N113").Value = -10681.8
$ws.Range("H126").Value = 6583
$ws.Range("I126").Value = 5998.5
$ws.Range("J126").Value = 6699.9
$ws.Range("K126").Value = 17995.5
$ws.Range("L126").Value = 20099.7
$ws.Range("M126").Value = -15525.5
$ws.Range("N126").Value = -25039.7
$ws.Range("H136").Value = 458231.38
$ws.Range("I136").Value = 834739.25
$ws.Range("K136").Value = 2504217.75
$ws.Range("M136").Value = -2501667.75

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 340246.5
$ws.Range("I128").Value = 340246.5
$ws.Range("K128").Value = 1020739.5
$ws.Range("M128").Value = -1015759.5
$ws.Range("H131").Value = 4724.5713
$ws.Range("I131").Value = 1732.3334
$ws.Range("K131").Value = 5197.0002
$ws.Range("M131").Value = -157.0002000000004

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7902.84
$ws.Range("I102").Value = 8636.190000000001
$ws.Range("J102").Value = 4052.75
$ws.Range("K102").Value = 8636.190000000001
$ws.Range("L102").Value = 4052.75
$ws.Range("M102").Value = -7014.190000000001
$ws.Range("N102").Value = -7296.75
$ws.Range("H126").Value = 5270.6665
$ws.Range("I126").Value = 2912
$ws.Range("J126").Value = 6450
$ws.Range("K126").Value = 8736
$ws.Range("L126").Value = 19350
$ws.Range("M126").Value = -6266
$ws.Range("N126").Value = -24290
$ws.Range("H132").Value = 775575.0600000001
$ws.Range("I132").Value = 1117164
$ws.Range("K132").Value = 3351492
$ws.Range("M132").Value = -3348962

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1188.6666
$ws.Range("I93").Value = 1126.5
$ws.Range("K93").Value = 1126.5
$ws.Range("M93").Value = 121.5
$ws.Range("H122").Value = 783063.1
$ws.Range("I122").Value = 718325.3
$ws.Range("J122").Value = 852780.9
$ws.Range("K122").Value = 2154975.9
$ws.Range("L122").Value = 2558342.7
$ws.Range("M122").Value = -2152525.9
$ws.Range("N122").Value = -2563242.7
$ws.Range("H132").Value = 5443.769
$ws.Range("J132").Value = 6375
$ws.Range("L132").Value = 19125
$ws.Range("N132").Value = -24185
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 33436.098
$ws.Range("I107").Value = 42878.168
$ws.Range("J107").Value = 1063.2858
$ws.Range("K107").Value = 128634.504
$ws.Range("L107").Value = 3189.8574
$ws.Range("M107").Value = -126714.504
$ws.Range("N107").Value = -7029.857400000001
$ws.Range("H132").Value = 2464.5107
$ws.Range("I132").Value = 1645.8
$ws.Range("K132").Value = 4937.4
$ws.Range("M132").Value = -2407.4
